$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 123456
$ws.Range("C3").Value = 566666
$ws.Range("E10").Select()
